$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 2.29992430780942
$ws.Range("D2").Value = 12.1030995033794
$ws.Range("E2").Value = 0.114069236054635

# Row 3
$ws.Range("C3").Value = 3.41562886864834
$ws.Range("D3").Value = 22.4278797671216
$ws.Range("E3").Value = 0.33762751190837

# Row 4
$ws.Range("C4").Value = 2.00159469135176
$ws.Range("D4").Value = 12.5313570240863
$ws.Range("E4").Value = 0.281405235805149

# Row 5
$ws.Range("C5").Value = 2.10542445744699
$ws.Range("D5").Value = 12.9404089404265
$ws.Range("E5").Value = 0.103572647553897

# Row 6
$ws.Range("C6").Value = 0.715366197432824
$ws.Range("D6").Value = 4.30091027820249
$ws.Range("E6").Value = 0.0412355967625848
$ws.Range("F6").Value = 0.0005
$ws.Range("G6").Value = 0.003

# Row 7
$ws.Range("C7").Value = 2.10312260265187
$ws.Range("D7").Value = 19.2198406390073
$ws.Range("E7").Value = 0.277663751629274

$wb.Save()
